$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.513.16'
$ws.Range("E2").Value = '  -0.19%  '
$ws.Range("D3").Value = '3.667.94'
$ws.Range("E3").Value = '  -1.00%  '
$ws.Range("E4").Value = '  -0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '623.50'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -7.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.88'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.92%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("E9").Value = '  -1.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.23'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +1.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.59%  '
$ws.Range("E12").Value = '  -2.60%  '
$ws.Range("D13").Value = '4.285.50'
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.31'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.82%  '
$ws.Range("D15").Value = '3.669.84'
$ws.Range("E15").Value = '  -1.48%  '
$ws.Range("D16").Value = '69.512.92'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.118'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.51'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '15.90'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.34'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '468.86'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.651'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.62%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '79.66'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '
$ws.Range("D24").Value = '3.815.39'
$ws.Range("E24").Value = '  -1.01%  '
$ws.Range("E25").Value = '  +0.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.15'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("E27").Value = '  -4.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -5.95%  '
$ws.Range("E29").Value = '  -3.19%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.66'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.10%  '
$ws.Range("E31").Value = '  +0.08%  '
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '26.63'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.96'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.58%  '
$ws.Range("B34").Value = 'NEARProtocol'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.39'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.02%  '
$ws.Range("B35").Value = 'RenzoRestakedETH'
$ws.Range("C35").Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range("D35").Value = '3.672.94'
$ws.Range("E35").Value = '  -0.61%  '
$ws.Range("E36").Value = '  -3.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.28'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '178.33'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.998'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.19%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.23'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.27%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.80'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.25%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0893'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.923'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.97%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '29.32'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +5.62%  '
$ws.Range("B46").Value = 'OKB'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '46.90'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.44%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.87'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.06%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.000265'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -5.18%  '
$ws.Range("E50").Value = '  -5.44%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.260'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.01%  '
